$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: "promotions" -> "brochures"
$ws.Name = "brochures"

# Add new row: hyperlink (URL text) in column B, brochure display name in column A
$url = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/offshore/region-3/aibt/AIBT_Region3_Q2_Brochure_1APR-30JUN22_VOL1.0.pdf"
$ws.Range("B2").Value = $url
$ws.Hyperlinks.Add($ws.Range("B2"), $url)
$ws.Range("A2").Value = "AIBT Region3 Q2 Brochure 1APR-30JUN 22_VOL1.0"

# Match final selection state left behind in the sheet
$ws.Range("A3").Select()
